$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 15-19 with revised ESO load data
$ws.Range("D15").Value = 4855
$ws.Range("E15").Value = 4614
$ws.Range("F15").Value = 4523
$ws.Range("G15").Value = 4522
$ws.Range("H15").Value = 4558
$ws.Range("I15").Value = 4697
$ws.Range("J15").Value = 5162
$ws.Range("K15").Value = 5669
$ws.Range("L15").Value = 6071
$ws.Range("M15").Value = 6184
$ws.Range("N15").Value = 6055
$ws.Range("O15").Value = 5956
$ws.Range("P15").Value = 5850
$ws.Range("Q15").Value = 5770
$ws.Range("R15").Value = 5663
$ws.Range("S15").Value = 5645
$ws.Range("T15").Value = 5662
$ws.Range("U15").Value = 5894
$ws.Range("V15").Value = 6441
$ws.Range("W15").Value = 6470
$ws.Range("X15").Value = 6252
$ws.Range("Y15").Value = 5977
$ws.Range("Z15").Value = 5724
$ws.Range("AA15").Value = 5336

$ws.Range("D16").Value = 4990
$ws.Range("E16").Value = 4735
$ws.Range("F16").Value = 4638
$ws.Range("G16").Value = 4637
$ws.Range("H16").Value = 4676
$ws.Range("I16").Value = 4823
$ws.Range("J16").Value = 5316
$ws.Range("K16").Value = 5853
$ws.Range("L16").Value = 6280
$ws.Range("M16").Value = 6400
$ws.Range("N16").Value = 6288
$ws.Range("O16").Value = 6203
$ws.Range("P16").Value = 6112
$ws.Range("Q16").Value = 6042
$ws.Range("R16").Value = 5951
$ws.Range("S16").Value = 5935
$ws.Range("T16").Value = 5943
$ws.Range("U16").Value = 6047
$ws.Range("V16").Value = 6292
$ws.Range("W16").Value = 6305
$ws.Range("X16").Value = 6087
$ws.Range("Y16").Value = 5811
$ws.Range("Z16").Value = 5558
$ws.Range("AA16").Value = 5171

$ws.Range("D17").Value = 4835
$ws.Range("E17").Value = 4602
$ws.Range("F17").Value = 4514
$ws.Range("G17").Value = 4513
$ws.Range("H17").Value = 4549
$ws.Range("I17").Value = 4682
$ws.Range("J17").Value = 5129
$ws.Range("K17").Value = 5616
$ws.Range("L17").Value = 6002
$ws.Range("M17").Value = 6111
$ws.Range("N17").Value = 5966
$ws.Range("O17").Value = 5855
$ws.Range("P17").Value = 5736
$ws.Range("Q17").Value = 5646
$ws.Range("R17").Value = 5527
$ws.Range("S17").Value = 5506
$ws.Range("T17").Value = 5520
$ws.Range("U17").Value = 5709
$ws.Range("V17").Value = 6155
$ws.Range("W17").Value = 6178
$ws.Range("X17").Value = 5996
$ws.Range("Y17").Value = 5766
$ws.Range("Z17").Value = 5554
$ws.Range("AA17").Value = 5230

$ws.Range("D18").Value = 4942
$ws.Range("E18").Value = 4747
$ws.Range("F18").Value = 4659
$ws.Range("G18").Value = 4628
$ws.Range("H18").Value = 4619
$ws.Range("I18").Value = 4675
$ws.Range("J18").Value = 4933
$ws.Range("K18").Value = 5403
$ws.Range("L18").Value = 5984
$ws.Range("M18").Value = 6157
$ws.Range("N18").Value = 6087
$ws.Range("O18").Value = 5930
$ws.Range("P18").Value = 5772
$ws.Range("Q18").Value = 5669
$ws.Range("R18").Value = 5531
$ws.Range("S18").Value = 5539
$ws.Range("T18").Value = 5693
$ws.Range("U18").Value = 6004
$ws.Range("V18").Value = 6335
$ws.Range("W18").Value = 6272
$ws.Range("X18").Value = 6054
$ws.Range("Y18").Value = 5827
$ws.Range("Z18").Value = 5714
$ws.Range("AA18").Value = 5512

$ws.Range("D19").Value = 5039
$ws.Range("E19").Value = 4885
$ws.Range("F19").Value = 4797
$ws.Range("G19").Value = 4741
$ws.Range("H19").Value = 4725
$ws.Range("I19").Value = 4769
$ws.Range("J19").Value = 4932
$ws.Range("K19").Value = 5261
$ws.Range("L19").Value = 5765
$ws.Range("M19").Value = 6154
$ws.Range("N19").Value = 6278
$ws.Range("O19").Value = 6282
$ws.Range("P19").Value = 6166
$ws.Range("Q19").Value = 6034
$ws.Range("R19").Value = 5803
$ws.Range("S19").Value = 5694
$ws.Range("T19").Value = 5929
$ws.Range("U19").Value = 6256
$ws.Range("V19").Value = 6711
$ws.Range("W19").Value = 6715
$ws.Range("X19").Value = 6536
$ws.Range("Y19").Value = 6293
$ws.Range("Z19").Value = 6081
$ws.Range("AA19").Value = 5765

# Add new row 20 (day 23) with ESO load data
$ws.Range("A20").Value = 2026
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 23
$ws.Range("D20").Value = 4881
$ws.Range("E20").Value = 4682
$ws.Range("F20").Value = 4603
$ws.Range("G20").Value = 4620
$ws.Range("H20").Value = 4658
$ws.Range("I20").Value = 4846
$ws.Range("J20").Value = 5441
$ws.Range("K20").Value = 6034
$ws.Range("L20").Value = 6433
$ws.Range("M20").Value = 6564
$ws.Range("N20").Value = 6384
$ws.Range("O20").Value = 6300
$ws.Range("P20").Value = 6170
$ws.Range("Q20").Value = 6086
$ws.Range("R20").Value = 5948
$ws.Range("S20").Value = 5893
$ws.Range("T20").Value = 5984
$ws.Range("U20").Value = 6257
$ws.Range("V20").Value = 6690
$ws.Range("W20").Value = 6719
$ws.Range("X20").Value = 6523
$ws.Range("Y20").Value = 6290
$ws.Range("Z20").Value = 6078
$ws.Range("AA20").Value = 5718
